$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update input values (feedback fix) - formulas throughout the sheet
# recalculate automatically from these inputs.
$ws.Range("B1").Value = 2
$ws.Range("B2").Value = 4
$ws.Range("B3").Value = 6

# Move active selection to H15 to match the saved view state
$ws.Range("H15").Select()
